$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.02
$ws.Range("J2").Value = 3
$ws.Range("O2").Value = 1.61
$ws.Range("P2").Value = 1.44
$ws.Range("Q2").Value = 2.7
$ws.Range("W2").Value = 1.98
$ws.Range("S3").Value = 1.66
$ws.Range("F4").Value = 2.9
$ws.Range("I4").Value = 2.76
$ws.Range("M4").Value = 1.08
$ws.Range("V4").Value = 1.57
$ws.Range("X4").Value = 970
$ws.Range("F5").Value = 5.1
$ws.Range("G5").Value = 5.8
$ws.Range("H5").Value = 1.98
$ws.Range("I5").Value = 2.08
$ws.Range("L5").Value = 1.49
$ws.Range("M5").Value = 1.13
$ws.Range("N5").Value = 2.56
$ws.Range("O5").Value = 1.56
$ws.Range("P5").Value = 1.52
$ws.Range("Q5").Value = 2.7
$ws.Range("R5").Value = 1.18
$ws.Range("S5").Value = 5.6
$ws.Range("T5").Value = 2.16
$ws.Range("V5").Value = 1.92
$ws.Range("W5").Value = 1.21
$ws.Range("Y5").Value = 6.6
$ws.Range("AA5").Value = 26
$ws.Range("AB5").Value = 13.5
$ws.Range("AE5").Value = 30
$ws.Range("AF5").Value = 38
$ws.Range("AG5").Value = 24
$ws.Range("AJ5").Value = 160
$ws.Range("AK5").Value = 110
$ws.Range("AL5").Value = 170
$ws.Range("AN5").Value = 180
$ws.Range("AO5").Value = 28
$ws.Range("F7").Value = 2.9
$ws.Range("G7").Value = 3.3
$ws.Range("H7").Value = 2.5
$ws.Range("I7").Value = 2.8
$ws.Range("L7").Value = 1.45
$ws.Range("M7").Value = 1.08
$ws.Range("P7").Value = 1.76
$ws.Range("Q7").Value = 2.1
$ws.Range("T7").Value = 1.82
$ws.Range("V7").Value = 1.55
$ws.Range("W7").Value = 1.43
$ws.Range("Z7").Value = 19.5
$ws.Range("AA7").Value = 42
$ws.Range("AE7").Value = 36
$ws.Range("AI7").Value = 55
$ws.Range("AJ7").Value = 60
$ws.Range("AK7").Value = 40
$ws.Range("AL7").Value = 55
$ws.Range("AN7").Value = 40
$ws.Range("AO7").Value = 30
$ws.Range("G8").Value = 16
$ws.Range("H8").Value = 1.33
$ws.Range("J8").Value = 4.8
$ws.Range("N8").Value = 3.5
$ws.Range("Q8").Value = 1.88
$ws.Range("U8").Value = 1.54
$ws.Range("G9").Value = 1.45
$ws.Range("L9").Value = 1.21
$ws.Range("R9").Value = 1.73
$ws.Range("F10").Value = 2.02
$ws.Range("G10").Value = 2.26
$ws.Range("H10").Value = 3.65
$ws.Range("I10").Value = 4.3
$ws.Range("N10").Value = 3.55
$ws.Range("Q10").Value = 1.92
$ws.Range("T10").Value = 1.76
$ws.Range("V10").Value = 1.3
$ws.Range("W10").Value = 1.79
$ws.Range("Z10").Value = 34
$ws.Range("AC10").Value = 9.4
$ws.Range("F11").Value = 5.4
$ws.Range("G11").Value = 7
$ws.Range("H11").Value = 1.74
$ws.Range("I11").Value = 1.89
$ws.Range("J11").Value = 3.4
$ws.Range("K11").Value = 3.85
$ws.Range("N11").Value = 2.92
$ws.Range("O11").Value = 1.45
$ws.Range("P11").Value = 1.62
$ws.Range("Q11").Value = 2.36
$ws.Range("S11").Value = 4.4
$ws.Range("T11").Value = 2.1
$ws.Range("U11").Value = 1.74
$ws.Range("V11").Value = 2.12
$ws.Range("W11").Value = 1.18
$ws.Range("Y11").Value = 7.2
$ws.Range("Z11").Value = 11
$ws.Range("AB11").Value = 16
$ws.Range("AC11").Value = 9.800000000000001
$ws.Range("AF11").Value = 46
$ws.Range("AG11").Value = 25
$ws.Range("AH11").Value = 30
$ws.Range("AJ11").Value = 220
$ws.Range("AK11").Value = 130
$ws.Range("AL11").Value = 140
$ws.Range("AM11").Value = 230
$ws.Range("AO11").Value = 18
$ws.Range("F12").Value = 3.65
$ws.Range("G12").Value = 4
$ws.Range("I12").Value = 2.22
$ws.Range("K12").Value = 3.75
$ws.Range("O12").Value = 1.33
$ws.Range("P12").Value = 1.89
$ws.Range("Q12").Value = 1.98
$ws.Range("S12").Value = 3.5
$ws.Range("V12").Value = 1.81
$ws.Range("W12").Value = 1.33
$ws.Range("H13").Value = 2.22
$ws.Range("I13").Value = 2.24
$ws.Range("N13").Value = 3.6
$ws.Range("O13").Value = 1.35
$ws.Range("P13").Value = 1.87
$ws.Range("Q13").Value = 2.02
$ws.Range("R13").Value = 1.33
$ws.Range("S13").Value = 3.6
$ws.Range("T13").Value = 1.83
$ws.Range("V13").Value = 1.8
$ws.Range("W13").Value = 1.37
$ws.Range("AA13").Value = 28
$ws.Range("AC13").Value = 7.8
$ws.Range("AE13").Value = 25
$ws.Range("AO13").Value = 19.5
$ws.Range("I14").Value = 1.98
$ws.Range("N14").Value = 3.6
$ws.Range("AF14").Value = 980
$ws.Range("L15").Value = 1.31
$ws.Range("Q15").Value = 1.62
$ws.Range("R15").Value = 1.56
$ws.Range("S15").Value = 2.58
$ws.Range("AB15").Value = 9.800000000000001
$ws.Range("AD15").Value = 1000
$ws.Range("AH15").Value = 32
$ws.Range("AN15").Value = 5.1
$ws.Range("Q16").Value = 1.86
$ws.Range("F18").Value = 6.6
$ws.Range("G18").Value = 7.4
$ws.Range("K18").Value = 5.2
$ws.Range("N18").Value = 5.2
$ws.Range("S18").Value = 2.64
$ws.Range("U18").Value = 2.08
$ws.Range("W18").Value = 1.16
$ws.Range("F19").Value = 2.3
$ws.Range("G19").Value = 2.4
$ws.Range("S19").Value = 3.9
$ws.Range("T19").Value = 1.83
$ws.Range("U19").Value = 2.02
$ws.Range("H20").Value = 8.800000000000001
$ws.Range("I20").Value = 10.5
$ws.Range("J20").Value = 4.6
$ws.Range("K20").Value = 5.1
$ws.Range("P20").Value = 2.04
$ws.Range("S20").Value = 3.2
$ws.Range("T20").Value = 2.12
$ws.Range("U20").Value = 1.8
$ws.Range("X20").Value = 19.5
$ws.Range("Y20").Value = 30
$ws.Range("Z20").Value = 85
$ws.Range("T23").Value = 1.87
$ws.Range("U23").Value = 1.92
$ws.Range("Y23").Value = 12
$ws.Range("F24").Value = 3.4
$ws.Range("I24").Value = 2.38
$ws.Range("V24").Value = 1.72
$ws.Range("W24").Value = 1.25
$ws.Range("F25").Value = 1.37
$ws.Range("I25").Value = 11.5
$ws.Range("K25").Value = 5.6
$ws.Range("P25").Value = 2.08
$ws.Range("Q25").Value = 1.78
$ws.Range("S25").Value = 3.05
$ws.Range("T25").Value = 2.16
$ws.Range("U25").Value = 1.71
$ws.Range("V25").Value = 1.09
$ws.Range("Y25").Value = 34
$ws.Range("AA25").Value = 510
$ws.Range("AE25").Value = 210
$ws.Range("AI25").Value = 180
$ws.Range("AM25").Value = 210
$ws.Range("AO25").Value = 330
$ws.Range("S26").Value = 3.95
$ws.Range("F27").Value = 1.63
$ws.Range("G27").Value = 1.77
$ws.Range("H27").Value = 5.1
$ws.Range("I27").Value = 8
$ws.Range("J27").Value = 3.6
$ws.Range("K27").Value = 4.8
$ws.Range("M27").Value = 1.07
$ws.Range("N27").Value = 3.55
$ws.Range("O27").Value = 1.33
$ws.Range("P27").Value = 1.9
$ws.Range("Q27").Value = 1.82
$ws.Range("R27").Value = 1.33
$ws.Range("T27").Value = 1.96
$ws.Range("U27").Value = 1.89
$ws.Range("W27").Value = 2.3
$ws.Range("AB27").Value = 10
$ws.Range("AF27").Value = 12
$ws.Range("AG27").Value = 12.5
$ws.Range("AJ27").Value = 21
$ws.Range("AM27").Value = 170
$ws.Range("AN27").Value = 13.5
